$wb = $excel.ActiveWorkbook

# --- Matriz_Económico ---
$ws = $wb.Worksheets.Item("Matriz_Económico")
$ws.Range("D2").Value = 0.3333333333333333
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.3333333333333333
$ws.Range("N2").Value = 7
$ws.Range("D3").Value = 0.3333333333333333
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.3333333333333333
$ws.Range("N3").Value = 7
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = 3
$ws.Range("E4").Value = 3
$ws.Range("J4").Value = 3
$ws.Range("K4").Value = 3
$ws.Range("D5").Value = 0.3333333333333333
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.3333333333333333
$ws.Range("N5").Value = 7
$ws.Range("B6").Value = 3
$ws.Range("C6").Value = 3
$ws.Range("E6").Value = 3
$ws.Range("J6").Value = 3
$ws.Range("K6").Value = 3
$ws.Range("B7").Value = 3
$ws.Range("C7").Value = 3
$ws.Range("E7").Value = 3
$ws.Range("J7").Value = 3
$ws.Range("K7").Value = 3
$ws.Range("D10").Value = 0.3333333333333333
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.3333333333333333
$ws.Range("N10").Value = 7
$ws.Range("D11").Value = 0.3333333333333333
$ws.Range("F11").Value = 0.3333333333333333
$ws.Range("G11").Value = 0.3333333333333333
$ws.Range("M11").Value = 0.3333333333333333
$ws.Range("B13").Value = 3
$ws.Range("C13").Value = 3
$ws.Range("E13").Value = 3
$ws.Range("J13").Value = 3
$ws.Range("K13").Value = 3
$ws.Range("B14").Value = 0.1428571428571428
$ws.Range("C14").Value = 0.1428571428571428
$ws.Range("E14").Value = 0.1428571428571428
$ws.Range("J14").Value = 0.1428571428571428

# --- Pesos_Locales_Económico ---
$ws = $wb.Worksheets.Item("Pesos_Locales_Económico")
$ws.Range("B2").Value = 0.06859393436079969
$ws.Range("B3").Value = 0.06859393436079969
$ws.Range("B4").Value = 0.1451001396860344
$ws.Range("B5").Value = 0.06859393436079965
$ws.Range("B6").Value = 0.1451001396860344
$ws.Range("B7").Value = 0.1451001396860344
$ws.Range("B8").Value = 0.01428744057464852
$ws.Range("B9").Value = 0.01428744057464852
$ws.Range("B10").Value = 0.06859393436079965
$ws.Range("B11").Value = 0.0697866565851307
$ws.Range("B12").Value = 0.008867907233478526
$ws.Range("B13").Value = 0.1451001396860344
$ws.Range("B14").Value = 0.009419377695460395
$ws.Range("B15").Value = 0.01428744057464852
$ws.Range("B16").Value = 0.01428744057464852

# --- Resultados ---
$ws = $wb.Worksheets.Item("Resultados")
$ws.Range("B2").Value = 0.04593215377490437
$ws.Range("B3").Value = 0.07143316379118989
$ws.Range("B4").Value = 0.05336464016097637
$ws.Range("B5").Value = 0.08714293726804158
$ws.Range("B6").Value = 0.0811177470439894
$ws.Range("B7").Value = 0.08284788187775356
$ws.Range("B8").Value = 0.08387555050917361
$ws.Range("B9").Value = 0.02548050271461082
$ws.Range("B10").Value = 0.05508831793560737
$ws.Range("B11").Value = 0.09065297775663407
$ws.Range("B12").Value = 0.09435133512397476
$ws.Range("B13").Value = 0.08097171984351949
$ws.Range("B14").Value = 0.06927477871710272
$ws.Range("B15").Value = 0.04449533620054171
$ws.Range("B16").Value = 0.03397095728198044

# --- Ranking_Alternativas ---
$ws = $wb.Worksheets.Item("Ranking_Alternativas")
$ws.Range("B2").Value = 0.09435133512397476
$ws.Range("B3").Value = 0.09065297775663407
$ws.Range("A4").Value = "Jean y Marie Thierry"
$ws.Range("B4").Value = 0.08714293726804158
$ws.Range("A5").Value = "Marcelo Mena"
$ws.Range("B5").Value = 0.08387555050917361
$ws.Range("A6").Value = "Las Cañas"
$ws.Range("B6").Value = 0.08284788187775356
$ws.Range("A7").Value = "Laguna Verde"
$ws.Range("B7").Value = 0.0811177470439894
$ws.Range("A8").Value = "Puertas Negras"
$ws.Range("B8").Value = 0.08097171984351949
$ws.Range("A9").Value = "Cordillera"
$ws.Range("B9").Value = 0.07143316379118989
$ws.Range("A10").Value = "Quebrada Verde"
$ws.Range("B10").Value = 0.06927477871710272
$ws.Range("A11").Value = "Placeres"
$ws.Range("B11").Value = 0.05508831793560737
$ws.Range("A12").Value = "Esperanza"
$ws.Range("B12").Value = 0.05336464016097637
$ws.Range("A13").Value = "Baron"
$ws.Range("B13").Value = 0.04593215377490437
$ws.Range("A14").Value = "Reina Isabel 2"
$ws.Range("B14").Value = 0.04449533620054171
$ws.Range("B15").Value = 0.03397095728198044
$ws.Range("B16").Value = 0.02548050271461082
